$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "60.773.99"
Set-TextValue "E2" "  -2.33%  "
Set-TextValue "D3" "2.908.65"
Set-TextValue "E3" "  -3.50%  "
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "586.92"
Set-TextValue "E5" "  -1.25%  "
Set-TextValue "E6" "  +0.01%  "
Set-TextValue "E7" "  +0.00%  "
Set-TextValue "D8" "0.505"
Set-TextValue "E8" "  -2.62%  "
Set-TextValue "D9" "2.907.91"
Set-TextValue "E9" "  -3.46%  "
Set-TextValue "D10" "6.71"
Set-TextValue "E10" "  +5.39%  "
Set-TextValue "D11" "0.144"
Set-TextValue "E11" "  -3.44%  "
Set-TextValue "D12" "0.449"
Set-TextValue "E12" "  -2.18%  "
Set-TextValue "D13" "0.0000225"
Set-TextValue "E13" "  -3.20%  "
Set-TextValue "D14" "34.17"
Set-TextValue "E14" "  -0.69%  "
Set-TextValue "E15" "  +0.33%  "
Set-TextValue "D16" "3.393.25"
Set-TextValue "E16" "  -3.51%  "
Set-TextValue "D17" "6.83"
Set-TextValue "E17" "  -2.37%  "
Set-TextValue "D18" "60.721.36"
Set-TextValue "E18" "  -2.47%  "
Set-TextValue "D19" "2.908.54"
Set-TextValue "E19" "  -3.50%  "
Set-TextValue "D20" "427.67"
Set-TextValue "E20" "  -4.31%  "
Set-TextValue "E21" "  -3.82%  "
Set-TextValue "D23" "7.12"
Set-TextValue "E23" "  -3.82%  "
Set-TextValue "B24" "RenderToken"
Set-TextValue "C24" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D24" "11.14"
Set-TextValue "E24" "  +2.44%  "
Set-TextValue "B25" "Litecoin"
Set-TextValue "C25" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D25" "80.65"
Set-TextValue "E25" "  -2.01%  "
Set-TextValue "D26" "2.23"
Set-TextValue "E26" "  -0.67%  "
Set-TextValue "D27" "11.87"
Set-TextValue "E27" "  -1.48%  "
Set-TextValue "E28" "  +0.07%  "
Set-TextValue "D29" "7.27"
Set-TextValue "E29" "  +1.32%  "
Set-TextValue "E30" "  -0.11%  "
Set-TextValue "E31" "  +2.71%  "
Set-TextValue "D32" "2.63"
Set-TextValue "E32" "  -2.85%  "
Set-TextValue "D33" "26.56"
Set-TextValue "E33" "  -3.31%  "
Set-TextValue "E34" "  -3.01%  "
Set-TextValue "D35" "0.0₃0840"
Set-TextValue "E35" "  -1.08%  "
Set-TextValue "E36" "  -1.77%  "
Set-TextValue "D37" "5.69"
Set-TextValue "E37" "  -2.51%  "
Set-TextValue "E38" "  -0.83%  "
Set-TextValue "B39" "OKB"
Set-TextValue "C39" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D39" "49.31"
Set-TextValue "E39" "  -1.70%  "
Set-TextValue "B40" "dogwifhat"
Set-TextValue "C40" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D40" "2.96"
Set-TextValue "E40" "  -0.14%  "
Set-TextValue "E41" "  -3.65%  "
Set-TextValue "D42" "0.122"
Set-TextValue "E42" "  -1.40%  "
Set-TextValue "D43" "0.293"
Set-TextValue "E43" "  +2.77%  "
Set-TextValue "D44" "41.70"
Set-TextValue "E44" "  +1.74%  "
Set-TextValue "D45" "0.0349"
Set-TextValue "E45" "  -1.01%  "
Set-TextValue "D46" "371.80"
Set-TextValue "E46" "  -5.65%  "
Set-TextValue "B47" "Maker"
Set-TextValue "C47" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D47" "2.657.41"
Set-TextValue "E47" "  -2.56%  "
Set-TextValue "B48" "Monero"
Set-TextValue "C48" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D48" "133.09"
Set-TextValue "E48" "  -1.31%  "
Set-TextValue "E49" "  -0.04%  "
Set-TextValue "D50" "25.25"
Set-TextValue "E50" "  +6.44%  "
Set-TextValue "E51" "  -1.08%  "
